# Update building block types in the Genome assembly template
$wb = $excel.ActiveWorkbook

# --- Sheet with template metadata: bump version number ---
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.1.8"

# --- Sheet with the annotation table: rename building block headers ---
$ws = $wb.Worksheets.Item("4COM04_GenomeAssembly")

$ws.Range("B1").Value = "Characteristic [BioSample Accession Number]"
$ws.Range("E1").Value = "Component [data filtering software]"
$ws.Range("N1").Value = "Component [next generation sequencing instrument model]"
$ws.Range("Q1").Value = "Component [sequence assembly algorithm]"
$ws.Range("W1").Value = "Characteristic [sequence assembly name]"
